# Insert a new data row at row 148 (pushing the existing rows 148-239 down
# to 149-240) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(148).Insert()

$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(148, 3).Value = "Ñuble"
$ws.Cells.Item(148, 4).Value = 44762
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112043
$ws.Cells.Item(148, 7).Value = "Pepino ensalada"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 100
$ws.Cells.Item(148, 11).Value = 19000
$ws.Cells.Item(148, 12).Value = 20000
$ws.Cells.Item(148, 13).Value = 19500
$ws.Cells.Item(148, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 325
$ws.Cells.Item(148, 17).Value = 60
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Preserve the date-formatted style ("s=2" in the original sheet) used by
# column D for every data row, matching the rest of the column.
$ws.Cells.Item(148, 4).NumberFormat = $ws.Cells.Item(149, 4).NumberFormat
